$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 65 (old rows 65-71 shift down to 67-73).
$ws.Rows("65:66").Insert()

# Populate the new row 65 with the new weekly price record.
$ws.Cells.Item(65, 1).Value = 7
$ws.Cells.Item(65, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(65, 3).Value = "Ñuble"
$ws.Cells.Item(65, 4).Value = 44785
$ws.Cells.Item(65, 5).Value = 16
$ws.Cells.Item(65, 6).Value = 100112040
$ws.Cells.Item(65, 7).Value = "Cilantro"
$ws.Cells.Item(65, 8).Value = "Sin especificar"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 2200
$ws.Cells.Item(65, 11).Value = 700
$ws.Cells.Item(65, 12).Value = 800
$ws.Cells.Item(65, 13).Value = 791
$ws.Cells.Item(65, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(65, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(65, 16).Value = 791
$ws.Cells.Item(65, 17).Value = 1
$ws.Cells.Item(65, 18).Value = "Hortaliza"

# Populate the new row 66 with the new weekly price record.
$ws.Cells.Item(66, 1).Value = 7
$ws.Cells.Item(66, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(66, 3).Value = "Ñuble"
$ws.Cells.Item(66, 4).Value = 44785
$ws.Cells.Item(66, 5).Value = 16
$ws.Cells.Item(66, 6).Value = 100112040
$ws.Cells.Item(66, 7).Value = "Cilantro"
$ws.Cells.Item(66, 8).Value = "Sin especificar"
$ws.Cells.Item(66, 9).Value = "Segunda"
$ws.Cells.Item(66, 10).Value = 250
$ws.Cells.Item(66, 11).Value = 600
$ws.Cells.Item(66, 12).Value = 600
$ws.Cells.Item(66, 13).Value = 600
$ws.Cells.Item(66, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(66, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(66, 16).Value = 600
$ws.Cells.Item(66, 17).Value = 1
$ws.Cells.Item(66, 18).Value = "Hortaliza"
